# Update countries & provincias Spain
#
# 1) Swap the display order of "El Salvador" / "Kirguistan" (rows 76/77)
# 2) Swap the display order of "Islas Malvinas" / "Groenlandia" (rows 209/210)
# 3) Update the "Datos actualizados..." timestamp in A1
# 4) Refresh the numeric stats for several country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Kirguistan now comes before El Salvador ---------------------------
$ws.Range("A76").Value = "Kirguistan"
$ws.Range("A77").Value = "El Salvador"

# --- 2) Groenlandia now comes before Islas Malvinas ------------------------
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- 3) Timestamp update ----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 08:54"

# --- 4) Numeric data refresh -------------------------------------------------
$updates = @{
    38  = @{ B = 50414; C = 807; D = 23119; E = 25989; F = 0;  G = 23; H = 1306 }
    47  = @{ B = 33594; C = 210; D = 20305; E = 12353; F = 0;  G = 16; H = 936  }
    76  = @{ B = 8486;  C = 345; D = 2983;  E = 5391;  F = 0;  G = 13; H = 112  }
    77  = @{ B = 8307;  C = 0;   D = 4955;  E = 3123;  F = 0;  G = 0;  H = 229  }
    97  = @{ B = 4210;  C = 5;   D = 2885;  E = 736  }
    142 = @{ B = 963;   C = 5;   D = 841;   E = 107  }
    192 = @{ B = 55;    C = 6;   D = 11;    E = 42   }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
